$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H3").Value = 34000
$ws_ALC.Range("J3").Value = 34000
$ws_ALC.Range("L3").Value = 34000
$ws_ALC.Range("N3").Value = -34228
$ws_ALC.Range("H9").Value = 97.46154
$ws_ALC.Range("I9").Value = 97.46154
$ws_ALC.Range("K9").Value = 97.46154
$ws_ALC.Range("M9").Value = 71.53846
$ws_ALC.Range("H87").Value = 63875
$ws_ALC.Range("J87").Value = 63875
$ws_ALC.Range("L87").Value = 63875
$ws_ALC.Range("N87").Value = -66371
$ws_ALC.Range("H90").Value = 63875
$ws_ALC.Range("J90").Value = 63875
$ws_ALC.Range("L90").Value = 191625
$ws_ALC.Range("N90").Value = -204105
$ws_ALC.Range("H101").Value = 210.75
$ws_ALC.Range("I101").Value = 210.75
$ws_ALC.Range("K101").Value = 632.25
$ws_ALC.Range("M101").Value = 989.75
$ws_ALC.Range("H102").Value = 34000
$ws_ALC.Range("J102").Value = 34000
$ws_ALC.Range("L102").Value = 34000
$ws_ALC.Range("N102").Value = -40490
$ws_ALC.Range("H111").Value = 275
$ws_ALC.Range("I111").Value = 275
$ws_ALC.Range("J111").Value = 0
$ws_ALC.Range("K111").Value = 825
$ws_ALC.Range("L111").Value = 0
$ws_ALC.Range("M111").Value = 2242
$ws_ALC.Range("N111").ClearContents()
$ws_ALC.Range("H113").Value = 4427.857
$ws_ALC.Range("I113").Value = 4299
$ws_ALC.Range("K113").Value = 4299
$ws_ALC.Range("M113").Value = -1045
$ws_ALC.Range("H121").Value = 5000
$ws_ALC.Range("J121").Value = 5000
$ws_ALC.Range("L121").Value = 15000
$ws_ALC.Range("N121").Value = -18494
$ws_ALC.Range("H125").Value = 4997.5
$ws_ALC.Range("I125").Value = 4995
$ws_ALC.Range("J125").Value = 5000
$ws_ALC.Range("K125").Value = 44955
$ws_ALC.Range("L125").Value = 45000
$ws_ALC.Range("M125").Value = -42495
$ws_ALC.Range("N125").Value = -49920
$ws_ALC.Range("H135").Value = 7573.304
$ws_ALC.Range("I135").Value = 3405.125
$ws_ALC.Range("J135").Value = 9796.333000000001
$ws_ALC.Range("K135").Value = 30646.125
$ws_ALC.Range("L135").Value = 88166.997
$ws_ALC.Range("M135").Value = -28111.125
$ws_ALC.Range("N135").Value = -93236.997
$ws_ALC.Range("H137").Value = 6005.4814
$ws_ALC.Range("I137").Value = 2390.3076
$ws_ALC.Range("J137").Value = 100000
$ws_ALC.Range("K137").Value = 7170.9228
$ws_ALC.Range("L137").Value = 300000
$ws_ALC.Range("M137").Value = -4620.9228
$ws_ALC.Range("N137").Value = -305100

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H45").Value = 4050
$ws_ARM.Range("I45").Value = 3000
$ws_ARM.Range("J45").Value = 4400
$ws_ARM.Range("K45").Value = 3000
$ws_ARM.Range("L45").Value = 4400
$ws_ARM.Range("M45").Value = -2623
$ws_ARM.Range("N45").Value = -5154
$ws_ARM.Range("H61").Value = 2232.2104
$ws_ARM.Range("I61").Value = 2200.75
$ws_ARM.Range("K61").Value = 2200.75
$ws_ARM.Range("M61").Value = -1988.75
$ws_ARM.Range("H74").Value = 5290.309
$ws_ARM.Range("I74").Value = 2978.6584
$ws_ARM.Range("J74").Value = 12060.143
$ws_ARM.Range("K74").Value = 2978.6584
$ws_ARM.Range("L74").Value = 12060.143
$ws_ARM.Range("M74").Value = -2104.6584
$ws_ARM.Range("N74").Value = -13808.143
$ws_ARM.Range("H77").Value = 5290.309
$ws_ARM.Range("I77").Value = 2978.6584
$ws_ARM.Range("J77").Value = 12060.143
$ws_ARM.Range("K77").Value = 14893.292
$ws_ARM.Range("L77").Value = 60300.715
$ws_ARM.Range("M77").Value = -10525.292
$ws_ARM.Range("N77").Value = -69036.715
$ws_ARM.Range("H102").Value = 1674.6666
$ws_ARM.Range("I102").Value = 1674.6666
$ws_ARM.Range("K102").Value = 1674.6666
$ws_ARM.Range("M102").Value = -52.66660000000002
$ws_ARM.Range("H132").Value = 7602.1113
$ws_ARM.Range("I132").Value = 5375.643
$ws_ARM.Range("K132").Value = 16126.929
$ws_ARM.Range("M132").Value = -13596.929
$ws_ARM.Range("H133").Value = 49900
$ws_ARM.Range("J133").Value = 49900
$ws_ARM.Range("L133").Value = 49900
$ws_ARM.Range("N133").Value = -54960
$ws_ARM.Range("H134").Value = 112599.2
$ws_ARM.Range("J134").Value = 112599.2
$ws_ARM.Range("L134").Value = 112599.2
$ws_ARM.Range("N134").Value = -122739.2
$ws_ARM.Range("H136").Value = 2232.2104
$ws_ARM.Range("I136").Value = 2200.75
$ws_ARM.Range("K136").Value = 6602.25
$ws_ARM.Range("M136").Value = -4052.25

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 4675.9
$ws_BSM.Range("I94").Value = 4527.5293
$ws_BSM.Range("J94").Value = 5516.6665
$ws_BSM.Range("K94").Value = 4527.5293
$ws_BSM.Range("L94").Value = 5516.6665
$ws_BSM.Range("M94").Value = -4076.5293
$ws_BSM.Range("N94").Value = -6418.6665
$ws_BSM.Range("H105").Value = 5748.2856
$ws_BSM.Range("I105").Value = 1553.9
$ws_BSM.Range("K105").Value = 1553.9
$ws_BSM.Range("M105").Value = 193.0999999999999
$ws_BSM.Range("H132").Value = 99198
$ws_BSM.Range("J132").Value = 99198
$ws_BSM.Range("L132").Value = 99198
$ws_BSM.Range("N132").Value = -109318
$ws_BSM.Range("H134").Value = 1880.4
$ws_BSM.Range("I134").Value = 1997
$ws_BSM.Range("J134").Value = 1414
$ws_BSM.Range("K134").Value = 5991
$ws_BSM.Range("L134").Value = 4242
$ws_BSM.Range("M134").Value = -3456
$ws_BSM.Range("N134").Value = -9312
$ws_BSM.Range("H135").Value = 68729
$ws_BSM.Range("J135").Value = 68729
$ws_BSM.Range("L135").Value = 68729
$ws_BSM.Range("N135").Value = -78869

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H86").Value = 94993.8
$ws_CRP.Range("J86").Value = 6841.6665
$ws_CRP.Range("L86").Value = 6841.6665
$ws_CRP.Range("N86").Value = -9087.666499999999
$ws_CRP.Range("H89").Value = 94993.8
$ws_CRP.Range("J89").Value = 6841.6665
$ws_CRP.Range("L89").Value = 34208.3325
$ws_CRP.Range("N89").Value = -45440.3325
$ws_CRP.Range("H134").Value = 2648.4167
$ws_CRP.Range("I134").Value = 1973.375
$ws_CRP.Range("J134").Value = 3998.5
$ws_CRP.Range("K134").Value = 5920.125
$ws_CRP.Range("L134").Value = 11995.5
$ws_CRP.Range("M134").Value = -3385.125
$ws_CRP.Range("N134").Value = -17065.5

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H9").Value = 1163745.1
$ws_CUL.Range("I9").Value = 950000
$ws_CUL.Range("J9").Value = 1200597.8
$ws_CUL.Range("K9").Value = 2850000
$ws_CUL.Range("L9").Value = 3601793.4
$ws_CUL.Range("M9").Value = -2849776
$ws_CUL.Range("N9").Value = -3602241.4

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H70").Value = 33538.08
$ws_GSM.Range("I70").Value = 40555.11
$ws_GSM.Range("J70").Value = 17749.75
$ws_GSM.Range("K70").Value = 40555.11
$ws_GSM.Range("L70").Value = 17749.75
$ws_GSM.Range("M70").Value = -40285.11
$ws_GSM.Range("N70").Value = -18289.75
$ws_GSM.Range("H73").Value = 33538.08
$ws_GSM.Range("I73").Value = 40555.11
$ws_GSM.Range("J73").Value = 17749.75
$ws_GSM.Range("K73").Value = 40555.11
$ws_GSM.Range("L73").Value = 17749.75
$ws_GSM.Range("M73").Value = -39619.11
$ws_GSM.Range("N73").Value = -19621.75
$ws_GSM.Range("H80").Value = 5851.6
$ws_GSM.Range("I80").Value = 5002.5
$ws_GSM.Range("K80").Value = 5002.5
$ws_GSM.Range("M80").Value = -4004.5
$ws_GSM.Range("H83").Value = 5851.6
$ws_GSM.Range("I83").Value = 5002.5
$ws_GSM.Range("K83").Value = 25012.5
$ws_GSM.Range("M83").Value = -20020.5
$ws_GSM.Range("H102").Value = 5185.2
$ws_GSM.Range("I102").Value = 6937.3335
$ws_GSM.Range("K102").Value = 6937.3335
$ws_GSM.Range("M102").Value = -5315.3335
$ws_GSM.Range("H125").Value = 50000
$ws_GSM.Range("J125").Value = 50000
$ws_GSM.Range("L125").Value = 50000
$ws_GSM.Range("N125").Value = -54920
$ws_GSM.Range("H126").Value = 2742.1428
$ws_GSM.Range("I126").Value = 2336.625
$ws_GSM.Range("J126").Value = 3282.8333
$ws_GSM.Range("K126").Value = 7009.875
$ws_GSM.Range("L126").Value = 9848.499899999999
$ws_GSM.Range("M126").Value = -4539.875
$ws_GSM.Range("N126").Value = -14788.4999
$ws_GSM.Range("H132").Value = 15054.593
$ws_GSM.Range("I132").Value = 15468.23
$ws_GSM.Range("J132").Value = 4300
$ws_GSM.Range("K132").Value = 46404.69
$ws_GSM.Range("L132").Value = 12900
$ws_GSM.Range("M132").Value = -43874.69
$ws_GSM.Range("N132").Value = -17960
$ws_GSM.Range("H134").Value = 48999.5
$ws_GSM.Range("J134").Value = 48999.5
$ws_GSM.Range("L134").Value = 146998.5
$ws_GSM.Range("N134").Value = -152068.5

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 2280.5386
$ws_LTW.Range("I7").Value = 1877.4546
$ws_LTW.Range("J7").Value = 4497.5
$ws_LTW.Range("K7").Value = 1877.4546
$ws_LTW.Range("L7").Value = 4497.5
$ws_LTW.Range("M7").Value = -1765.4546
$ws_LTW.Range("N7").Value = -4721.5
$ws_LTW.Range("H40").Value = 2739.2
$ws_LTW.Range("I40").Value = 2710.2222
$ws_LTW.Range("J40").Value = 3000
$ws_LTW.Range("K40").Value = 2710.2222
$ws_LTW.Range("L40").Value = 3000
$ws_LTW.Range("M40").Value = -2574.2222
$ws_LTW.Range("N40").Value = -3272
$ws_LTW.Range("H82").Value = 2113.75
$ws_LTW.Range("I82").Value = 2169.842
$ws_LTW.Range("J82").Value = 1900.6
$ws_LTW.Range("K82").Value = 2169.842
$ws_LTW.Range("L82").Value = 1900.6
$ws_LTW.Range("M82").Value = -1808.842
$ws_LTW.Range("N82").Value = -2622.6
$ws_LTW.Range("H85").Value = 2113.75
$ws_LTW.Range("I85").Value = 2169.842
$ws_LTW.Range("J85").Value = 1900.6
$ws_LTW.Range("K85").Value = 2169.842
$ws_LTW.Range("L85").Value = 1900.6
$ws_LTW.Range("M85").Value = -921.8420000000001
$ws_LTW.Range("N85").Value = -4396.6
$ws_LTW.Range("H93").Value = 3346.6
$ws_LTW.Range("I93").Value = 1355.5714
$ws_LTW.Range("K93").Value = 1355.5714
$ws_LTW.Range("M93").Value = -107.5714
$ws_LTW.Range("H122").Value = 3796
$ws_LTW.Range("I122").Value = 3194.5
$ws_LTW.Range("J122").Value = 4999
$ws_LTW.Range("K122").Value = 9583.5
$ws_LTW.Range("L122").Value = 14997
$ws_LTW.Range("M122").Value = -7133.5
$ws_LTW.Range("N122").Value = -19897
$ws_LTW.Range("H126").Value = 2280.5386
$ws_LTW.Range("I126").Value = 1877.4546
$ws_LTW.Range("J126").Value = 4497.5
$ws_LTW.Range("K126").Value = 5632.3638
$ws_LTW.Range("L126").Value = 13492.5
$ws_LTW.Range("M126").Value = -3162.3638
$ws_LTW.Range("N126").Value = -18432.5
$ws_LTW.Range("H132").Value = 2392.75
$ws_LTW.Range("I132").Value = 1624.5
$ws_LTW.Range("J132").Value = 3380.5
$ws_LTW.Range("K132").Value = 4873.5
$ws_LTW.Range("L132").Value = 10141.5
$ws_LTW.Range("M132").Value = -2343.5
$ws_LTW.Range("N132").Value = -15201.5
$ws_LTW.Range("H136").Value = 11009.25
$ws_LTW.Range("I136").Value = 4910
$ws_LTW.Range("J136").Value = 14058.875
$ws_LTW.Range("K136").Value = 14730
$ws_LTW.Range("L136").Value = 42176.625
$ws_LTW.Range("M136").Value = -12180
$ws_LTW.Range("N136").Value = -47276.625

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H107").Value = 1254.6511
$ws_WVR.Range("I107").Value = 877.9231
$ws_WVR.Range("K107").Value = 2633.7693
$ws_WVR.Range("M107").Value = -713.7692999999999
$ws_WVR.Range("H132").Value = 2728.7144
$ws_WVR.Range("I132").Value = 2359.4443
$ws_WVR.Range("K132").Value = 7078.3329
$ws_WVR.Range("M132").Value = -4548.3329

Write-Output "Applied 269 cell updates across 8 sheets"